$wb = $excel.ActiveWorkbook

# LCFS updates for FF55: switch to the PVTStL sheet (becomes the active tab)
$ws = $wb.Worksheets.Item("PVTStL")
$ws.Activate()

# PVsTL: mark "ships" (row 6) as a vehicle type subject to LCFS
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1

# Reflect the cursor/selection position left on the PVTStL sheet
$ws.Range("B6").Select()
